$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 2938.103010863317
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 74466.72682088954
